$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.334.84"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.648.30"
$ws.Range("E3").Value = "  +2.62%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.40"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.93"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.647.41"
$ws.Range("E9").Value = "  +2.63%  "

$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.41"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.123.57"
$ws.Range("E15").Value = "  +2.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.187.70"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.634.15"
$ws.Range("E18").Value = "  +1.98%  "

$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  +3.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.58"
$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.83"
$ws.Range("E22").Value = "  +2.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.03"
$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("E25").Value = "  +2.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.53"
$ws.Range("E26").Value = "  -3.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.71"
$ws.Range("E27").Value = "  +5.45%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "543.02"
$ws.Range("E29").Value = "  +15.24%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  -0.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  +4.61%  "

$ws.Range("E33").Value = "  +7.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0813"
$ws.Range("E34").Value = "  +1.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "171.74"
$ws.Range("E35").Value = "  -2.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("E36").Value = "  +14.22%  "

$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.13"
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  +6.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.10"
$ws.Range("E41").Value = "  +6.25%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.45"
$ws.Range("E44").Value = "  +3.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0580"
$ws.Range("E45").Value = "  +7.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.632"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.82"
$ws.Range("E49").Value = "  +3.15%  "

$ws.Range("E50").Value = "  +2.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.23"
$ws.Range("E51").Value = "  -1.21%  "
